$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form Responses 1")

# ---------------------------------------------------------------------------
# 1. Workbook window width change
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Width = 20980

# ---------------------------------------------------------------------------
# 2. Append two new directory rows (167 and 168)
# ---------------------------------------------------------------------------

# ---- Row 167 : Rue Sparks / Daylight Chasers -------------------------------
$ws.Cells.Item(167, 1).Value = 44253.887003136573
$ws.Cells.Item(167, 2).Value = "Queer Lit"
$ws.Cells.Item(167, 3).Value = "Rue Sparks"
$ws.Cells.Item(167, 4).Value = "Speculative, Magical Realism, SFF"
$ws.Cells.Item(167, 5).Value = "Daylight Chasers"
$ws.Cells.Item(167, 6).Value = "https://linktr.ee/ruesparks"
$ws.Cells.Item(167, 7).Value = "During a time-zone-hopping road trip, a client's mercurial moods and thinly veiled secret leaves Keenan wondering: how can he be the guide when even he is feeling lost?"
$ws.Cells.Item(167, 8).Value = "@sparks_writes"
$ws.Cells.Item(167, 9).Formula = '="["&CHAR(39)&C167&CHAR(39)&","&CHAR(39)&F167&CHAR(39)&","&CHAR(39)&E167&CHAR(39)&","&CHAR(39)&D167&CHAR(39)&","&CHAR(39)&H167&CHAR(39)&","&CHAR(39)&G167&CHAR(39)&","&CHAR(39)&J167&CHAR(39)&","&CHAR(39)&K167&CHAR(39)&"],"'
$ws.Cells.Item(167, 10).Value = "genre-speculative"
$ws.Cells.Item(167, 11).Value = "age-adult"

# ---- Row 168 : Miles Nelson / Riftmaster ------------------------------------
$ws.Cells.Item(168, 1).Value = 44264.16777443287
$ws.Cells.Item(168, 2).Value = "Queer Lit"
$ws.Cells.Item(168, 3).Value = "Miles Nelson"
$ws.Cells.Item(168, 4).Value = "Science Fiction"
$ws.Cells.Item(168, 5).Value = "Riftmaster"
$ws.Cells.Item(168, 6).Value = "https://www.amazon.com/Riftmaster-Miles-Nelson-ebook/dp/B08WJGPY3W/"
$ws.Cells.Item(168, 7).Value = "College student bailey jones is wrenched away from earth by a mysterious and unpredictable force known as the Rift. While stranded on an alien planet, he meets a mysterious traveller known as the Riftmaster."
$ws.Cells.Item(168, 8).Value = "@ProbablyMiles"
$ws.Cells.Item(168, 9).Formula = '="["&CHAR(39)&C168&CHAR(39)&","&CHAR(39)&F168&CHAR(39)&","&CHAR(39)&E168&CHAR(39)&","&CHAR(39)&D168&CHAR(39)&","&CHAR(39)&H168&CHAR(39)&","&CHAR(39)&G168&CHAR(39)&","&CHAR(39)&J168&CHAR(39)&","&CHAR(39)&K168&CHAR(39)&"],"'
$ws.Cells.Item(168, 10).Value = "genre-speculative"
$ws.Cells.Item(168, 11).Value = "age-adult"

# ---------------------------------------------------------------------------
# 3. Copy cell formatting down from the previous last row (166) so the new
#    rows are styled the same way as the rest of the directory.
# ---------------------------------------------------------------------------
for ($col = 1; $col -le 11; $col++) {
    $ws.Cells.Item(166, $col).Copy()
    $ws.Cells.Item(167, $col).PasteSpecial(-4122)
    $ws.Cells.Item(168, $col).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Hyperlinks for the new "Link to Book(s)" cells
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Cells.Item(167, 6), "https://linktr.ee/ruesparks") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(168, 6), "https://www.amazon.com/Riftmaster-Miles-Nelson-ebook/dp/B08WJGPY3W/") | Out-Null

# ---------------------------------------------------------------------------
# 5. Sheet view tweaks (frozen pane scroll position & selected cell)
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A143").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("I169").Select()
